$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("I3").Value = -0.1390649750100463
$ws.Range("J3").Value = 0.6440599939930639
$ws.Range("K3").Value = 0.4936289468824919
$ws.Range("L3").Value = 2.423707911152657

$ws.Range("I20").Value = 0.1070944970656949
$ws.Range("J20").Value = 0.6482476954050463
$ws.Range("K20").Value = 0.1388475111057705
$ws.Range("L20").Value = 2.229431919465588
